$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the DEPAKINE CHRONO row (row 36): H (ratio), L (price) and N (count) ---
$ws.Cells.Item(36, 8).Value = "1:0"
$ws.Cells.Item(36, 12).Value = 139.68
$ws.Cells.Item(36, 14).Value = 1

# --- Remove the "سلاكه اسنان بلاستك" line item (row 111) entirely, shifting the ---
# --- rows below it up by one (merged cells / shared strings follow automatically) ---
$ws.Rows.Item(111).Delete()

# --- Recalculate the displayed total (now on row 117 after the shift) ---
$ws.Cells.Item(117, 11).Value = 6356

# --- The totals row and the footer row recompute their auto-fit heights once the ---
# --- item row above them is removed ---
$ws.Rows.Item(117).RowHeight = 26.25
$ws.Rows.Item(118).RowHeight = 16.5
